# Apply hybrid bold + color (2C3E50) highlighting to quantitative metrics
# across specific bullet points in the resume, per the target diff.
#
# Strategy: for each target paragraph, locate it by matching its leading
# text (stable even as earlier edits shift character offsets), then for
# each metric substring inside that paragraph (in left-to-right order),
# re-fetch a fresh Range over that paragraph and use Find.Execute to
# narrow the range down to just that substring, then apply Bold + Color.
# Word's Find.Execute collapses/narrows the range to the match itself,
# exactly like Word does interactively.

$d = $word.ActiveDocument

$highlightColor = 5258796   # 0x503E2C == BGR(0x2C,0x3E,0x50) -> w:color "2C3E50"

function Get-ParagraphByLeadingText($doc, [string]$needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs($i)
        if ($p.Range.Text -like "*$needle*") {
            return $p
        }
    }
    return $null
}

function Highlight-Metric($doc, $paragraph, [string]$metricText) {
    $rng = $paragraph.Range
    $found = $rng.Find.Execute($metricText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Font.Bold = $true
        $rng.Font.Color = $highlightColor
    }
    return $found
}

# 1) "Discovered systematic race coding errors..." -> 23%, 64%
$p = Get-ParagraphByLeadingText $d "Discovered systematic race coding errors"
Highlight-Metric $d $p "23%" | Out-Null
Highlight-Metric $d $p "64%" | Out-Null

# 2) "Utilized advanced sampling methods..." -> ±4.2%, ±2.1%, 71%, 87%
$p = Get-ParagraphByLeadingText $d "Utilized advanced sampling methods"
Highlight-Metric $d $p "±4.2%" | Out-Null
Highlight-Metric $d $p "±2.1%" | Out-Null
Highlight-Metric $d $p "71%" | Out-Null
Highlight-Metric $d $p "87%" | Out-Null

# 3) "Trigonometric algorithm for boundary estimation..." -> 73.5%, $4.7M
$p = Get-ParagraphByLeadingText $d "Trigonometric algorithm for boundary estimation"
Highlight-Metric $d $p "73.5%" | Out-Null
Highlight-Metric $d $p "$4.7M" | Out-Null

# 4) "Built real-time FEC analysis systems..." -> $2 (of "$2 trillion")
$p = Get-ParagraphByLeadingText $d "Built real-time FEC analysis systems"
Highlight-Metric $d $p "$2" | Out-Null

# 5) "Modernized legacy ETL processes..." -> 57%
$p = Get-ParagraphByLeadingText $d "Modernized legacy ETL processes"
Highlight-Metric $d $p "57%" | Out-Null

# 6) "Revenue generation: Delivered $4.9M..." -> $4.9M
$p = Get-ParagraphByLeadingText $d "Revenue generation: Delivered"
Highlight-Metric $d $p "$4.9M" | Out-Null

# 7) "23% conversion rate improvement" -> 23%
$p = Get-ParagraphByLeadingText $d "conversion rate improvement"
Highlight-Metric $d $p "23%" | Out-Null

# 8) "Platform impact: Built redistricting system serving 12,847..." -> 12,847
$p = Get-ParagraphByLeadingText $d "Platform impact: Built redistricting system serving"
Highlight-Metric $d $p "12,847" | Out-Null

Write-Output "Highlighting complete"
